# Auto-generated edit script
# Applies numeric cell updates to restore the recalculated leve-profit values
# across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR, matching the target diff.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("H32").Value = 1754.6364
$ws.Range("I32").Value = 1151
$ws.Range("J32").Value = 1888.7778
$ws.Range("K32").Value = 1151
$ws.Range("L32").Value = 1888.7778
$ws.Range("M32").Value = -825
$ws.Range("N32").Value = -2540.7778
$ws.Range("H98").Value = 2535.7693
$ws.Range("I98").Value = 1913.25
$ws.Range("J98").Value = 10006
$ws.Range("K98").Value = 1913.25
$ws.Range("L98").Value = 10006
$ws.Range("M98").Value = -415.25
$ws.Range("N98").Value = -13002
$ws.Range("H122").Value = 2535.7693
$ws.Range("I122").Value = 1913.25
$ws.Range("J122").Value = 10006
$ws.Range("K122").Value = 5739.75
$ws.Range("L122").Value = 30018
$ws.Range("M122").Value = -3289.75
$ws.Range("N122").Value = -34918
$ws.Range("H127").Value = 2878.9546
$ws.Range("I127").Value = 1848.5
$ws.Range("J127").Value = 2982
$ws.Range("K127").Value = 5545.5
$ws.Range("L127").Value = 8946
$ws.Range("M127").Value = -585.5
$ws.Range("N127").Value = -18866
$ws.Range("H132").Value = 1687.8572
$ws.Range("I132").Value = 1714.2122
$ws.Range("K132").Value = 5142.6366
$ws.Range("M132").Value = -2612.6366
$ws.Range("H138").Value = 4252.69
$ws.Range("I138").Value = 1323
$ws.Range("J138").Value = 4810.726
$ws.Range("K138").Value = 3969
$ws.Range("L138").Value = 14432.178
$ws.Range("M138").Value = 1171
$ws.Range("N138").Value = -24712.178

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("H57").Value = 7885.7144
$ws.Range("I57").Value = 7885.7144
$ws.Range("K57").Value = 7885.7144
$ws.Range("M57").Value = -7401.7144
$ws.Range("H74").Value = 810.6923
$ws.Range("I74").Value = 853.6
$ws.Range("J74").Value = 752.1818
$ws.Range("K74").Value = 853.6
$ws.Range("L74").Value = 752.1818
$ws.Range("M74").Value = 20.39999999999998
$ws.Range("N74").Value = -2500.1818
$ws.Range("H77").Value = 810.6923
$ws.Range("I77").Value = 853.6
$ws.Range("J77").Value = 752.1818
$ws.Range("K77").Value = 4268
$ws.Range("L77").Value = 3760.909
$ws.Range("M77").Value = 100
$ws.Range("N77").Value = -12496.909
$ws.Range("H122").Value = 2500
$ws.Range("I122").Value = 2000
$ws.Range("K122").Value = 6000
$ws.Range("M122").Value = -3550

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("H94").Value = 1312.9
$ws.Range("I94").Value = 1099.1428
$ws.Range("K94").Value = 1099.1428
$ws.Range("M94").Value = -648.1428000000001
$ws.Range("H97").Value = 9400
$ws.Range("I97").Value = 9400
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 9400
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()
$ws.Range("H134").Value = 21729.32
$ws.Range("I134").Value = 1572.6154
$ws.Range("J134").Value = 93194
$ws.Range("K134").Value = 4717.8462
$ws.Range("L134").Value = 279582
$ws.Range("M134").Value = -2182.8462
$ws.Range("N134").Value = -284652

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("H31").Value = 1399.1039
$ws.Range("I31").Value = 1089.7858
$ws.Range("J31").Value = 1770.2858
$ws.Range("K31").Value = 1089.7858
$ws.Range("L31").Value = 1770.2858
$ws.Range("M31").Value = -794.7858000000001
$ws.Range("N31").Value = -2360.2858
$ws.Range("H34").Value = 1399.1039
$ws.Range("I34").Value = 1089.7858
$ws.Range("J34").Value = 1770.2858
$ws.Range("K34").Value = 1089.7858
$ws.Range("L34").Value = 1770.2858
$ws.Range("M34").Value = -887.7858000000001
$ws.Range("N34").Value = -2174.2858
$ws.Range("H58").Value = 4076.7273
$ws.Range("I58").Value = 968.6
$ws.Range("J58").Value = 13789.625
$ws.Range("K58").Value = 968.6
$ws.Range("L58").Value = 13789.625
$ws.Range("M58").Value = -765.6
$ws.Range("N58").Value = -14195.625
$ws.Range("H60").Value = 18350
$ws.Range("J60").Value = 18350
$ws.Range("L60").Value = 18350
$ws.Range("N60").Value = -19372
$ws.Range("H99").Value = 2136.842
$ws.Range("I99").Value = 2082.353
$ws.Range("J99").Value = 2600
$ws.Range("K99").Value = 2082.353
$ws.Range("L99").Value = 2600
$ws.Range("M99").Value = -584.3530000000001
$ws.Range("N99").Value = -5596
$ws.Range("H122").Value = 834398.0600000001
$ws.Range("I122").Value = 1112118.1
$ws.Range("J122").Value = 1238
$ws.Range("K122").Value = 3336354.3
$ws.Range("L122").Value = 3714
$ws.Range("M122").Value = -3333904.3
$ws.Range("N122").Value = -8614
$ws.Range("H126").Value = 2136.842
$ws.Range("I126").Value = 2082.353
$ws.Range("J126").Value = 2600
$ws.Range("K126").Value = 6247.059
$ws.Range("L126").Value = 7800
$ws.Range("M126").Value = -3777.059
$ws.Range("N126").Value = -12740
$ws.Range("H136").Value = 4076.7273
$ws.Range("I136").Value = 968.6
$ws.Range("J136").Value = 13789.625
$ws.Range("K136").Value = 2905.8
$ws.Range("L136").Value = 41368.875
$ws.Range("M136").Value = -355.8000000000002
$ws.Range("N136").Value = -46468.875

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("H122").Value = 854.4545000000001
$ws.Range("I122").Value = 544.44446
$ws.Range("J122").Value = 2249.5
$ws.Range("K122").Value = 4900.00014
$ws.Range("L122").Value = 20245.5
$ws.Range("M122").Value = -2450.00014
$ws.Range("N122").Value = -25145.5
$ws.Range("H131").Value = 22286.102
$ws.Range("J131").Value = 1869.5
$ws.Range("L131").Value = 5608.5
$ws.Range("N131").Value = -15688.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("H94").Value = 30448
$ws.Range("J94").Value = 30448
$ws.Range("L94").Value = 30448
$ws.Range("N94").Value = -31800
$ws.Range("H122").Value = 2167333.5
$ws.Range("I122").Value = 3249000.2
$ws.Range("K122").Value = 9747000.600000001
$ws.Range("M122").Value = -9744550.600000001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("H40").Value = 685127.6
$ws.Range("I40").Value = 1026001.4
$ws.Range("J40").Value = 3380
$ws.Range("K40").Value = 1026001.4
$ws.Range("L40").Value = 3380
$ws.Range("M40").Value = -1025865.4
$ws.Range("N40").Value = -3652
$ws.Range("H93").Value = 2358.4
$ws.Range("J93").Value = 2446
$ws.Range("L93").Value = 2446
$ws.Range("N93").Value = -4942
$ws.Range("H119").Value = 36666.668
$ws.Range("J119").Value = 36666.668
$ws.Range("L119").Value = 36666.668
$ws.Range("N119").Value = -46342.668
$ws.Range("H132").Value = 2527706.5
$ws.Range("I132").Value = 3970416
$ws.Range("J132").Value = 2964.8333
$ws.Range("K132").Value = 11911248
$ws.Range("L132").Value = 8894.499899999999
$ws.Range("M132").Value = -11908718
$ws.Range("N132").Value = -13954.4999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item(8)
$ws.Range("H126").Value = 988.5
$ws.Range("I126").Value = 1021.6
$ws.Range("J126").Value = 933.3333
$ws.Range("K126").Value = 3064.8
$ws.Range("L126").Value = 2799.9999
$ws.Range("M126").Value = -594.8000000000002
$ws.Range("N126").Value = -7739.9999
$ws.Range("H132").Value = 1831.7778
$ws.Range("I132").Value = 792.4
$ws.Range("J132").Value = 3131
$ws.Range("K132").Value = 2377.2
$ws.Range("L132").Value = 9393
$ws.Range("M132").Value = 152.8000000000002
$ws.Range("N132").Value = -14453
$ws.Range("H136").Value = 1493.7188
$ws.Range("I136").Value = 1351.2632
$ws.Range("J136").Value = 1701.9231
$ws.Range("K136").Value = 4053.7896
$ws.Range("L136").Value = 5105.7693
$ws.Range("M136").Value = -1503.7896
$ws.Range("N136").Value = -10205.7693
